# Updated remaining queries for C3DC
# - Rewrites the JOIN conditions in all SQL queries stored on Sheet1
#   (std.id / prt.id -> std.study_id / prt.participant_id, matching
#   renamed join columns on the source dataframes).
# - Widens column C slightly now that the bestFit measurement is stale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query([string]$cellRef) {
    $text = $ws.Range($cellRef).Text

    $text = $text.Replace(
        'df_participant prt ON std.id = prt."study.id"',
        'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace(
        'df_diagnoses dgn ON prt.id = dgn."participant.id"',
        'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace(
        'df_treatments trt ON prt.id = trt."participant.id"',
        'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace(
        'df_treatment_resp trr ON prt.id = trr."participant.id"',
        'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace(
        'df_survival srv ON prt.id = srv."participant.id"',
        'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace(
        'df_reference_files rfs ON std.id = rfs."study.id"',
        'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $ws.Range($cellRef).Value = $text
}

# C2 = StatQuery; B2:B7 = TabQuery for each tab row.
Update-Query("C2")
Update-Query("B2")
Update-Query("B3")
Update-Query("B4")
Update-Query("B5")
Update-Query("B6")
Update-Query("B7")

# Column C widened (bestFit measurement no longer applies after the text edits).
$ws.Columns.Item(3).ColumnWidth = 62.8333333
